# Data source corrected and updated
# - Columns J and K on Sheet1 were re-derived from the corrected source data:
#   J (previously header "r" + 0.3 filler) -> constant 0.5
#   K (previously header "s" + 0.5 filler) -> constant 1
# - View state (zoom / scroll position / selection) updated to match the
#   author's last on-screen state when the workbook was saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 51

# Column J: now a uniform 0.5 for every data row (header string "r" is gone).
$ws.Range("J1:J$lastRow").Value = 0.5

# Column K: now a uniform 1 for every data row (header string "s" is gone).
$ws.Range("K1:K$lastRow").Value = 1

# Restore the view state recorded in the saved workbook: zoomed to 90%,
# scrolled so row 19 is at the top, with K1:K51 selected (active cell K1).
$win = $excel.ActiveWindow
$win.ScrollRow = 19
$win.ScrollColumn = 1
$win.Zoom = 90

$ws.Range("K1:K$lastRow").Select()
